$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.482.94"
Set-TextValue "E2" "  +3.21%  "
Set-TextValue "D3" "1.819.33"
Set-TextValue "E3" "  +4.64%  "
Set-TextValue "E4" "  +0.75%  "
Set-TextValue "D5" "343.52"
Set-TextValue "E5" "  +2.60%  "
Set-TextValue "E6" "  +0.63%  "
Set-TextValue "D7" "0.3846"
Set-TextValue "E7" "  +3.28%  "
Set-TextValue "E8" "  +4.03%  "
Set-TextValue "D9" "49.00"
Set-TextValue "E9" "  -1.40%  "
Set-TextValue "D10" "1.236"
Set-TextValue "E10" "  +2.39%  "
Set-TextValue "D11" "0.07788"
Set-TextValue "E11" "  +3.28%  "
Set-TextValue "D12" "1.001"
Set-TextValue "E12" "  +0.84%  "
Set-TextValue "D13" "22.38"
Set-TextValue "E13" "  +8.77%  "
Set-TextValue "D14" "6.608"
Set-TextValue "E14" "  +2.61%  "
Set-TextValue "D15" "1.821.31"
Set-TextValue "E15" "  +4.99%  "
Set-TextValue "D16" "7.223"
Set-TextValue "E16" "  +2.88%  "
Set-TextValue "D17" "0.00001122"
Set-TextValue "E17" "  +1.94%  "
Set-TextValue "D18" "0.06717"
Set-TextValue "E18" "  +0.61%  "
Set-TextValue "D19" "86.46"
Set-TextValue "E19" "  +3.39%  "
Set-TextValue "D20" "1.0000"
Set-TextValue "E20" "  +0.65%  "
Set-TextValue "D21" "17.68"
Set-TextValue "E21" "  +5.14%  "
Set-TextValue "D22" "6.576"
Set-TextValue "E22" "  +6.33%  "
Set-TextValue "D23" "13.22"
Set-TextValue "E23" "  +0.50%  "
Set-TextValue "D24" "27.495.07"
Set-TextValue "E24" "  +3.49%  "
Set-TextValue "D25" "2.466"
Set-TextValue "E25" "  +0.02%  "
Set-TextValue "D26" "2.693"
Set-TextValue "E26" "  +6.62%  "
Set-TextValue "E27" "  +13.74%  "
Set-TextValue "D28" "1.473"
Set-TextValue "E28" "  +3.15%  "
Set-TextValue "D29" "154.02"
Set-TextValue "E29" "  +1.32%  "
Set-TextValue "D30" "2.023.83"
Set-TextValue "E30" "  +5.00%  "
Set-TextValue "D31" "136.52"
Set-TextValue "E31" "  +3.50%  "
Set-TextValue "D32" "6.389"
Set-TextValue "E32" "  +2.54%  "
Set-TextValue "D33" "4.064"
Set-TextValue "E33" "  -1.32%  "
Set-TextValue "D34" "13.95"
Set-TextValue "E34" "  +5.30%  "
Set-TextValue "D35" "0.08817"
Set-TextValue "E35" "  +2.58%  "
Set-TextValue "E36" "  -0.79%  "
Set-TextValue "D37" "5.625"
Set-TextValue "E37" "  +2.74%  "
Set-TextValue "D38" "0.7057"
Set-TextValue "E38" "  +12.55%  "
Set-TextValue "D39" "0.2267"
Set-TextValue "E39" "  +4.14%  "
Set-TextValue "E40" "  +2.16%  "
Set-TextValue "D41" "0.06485"
Set-TextValue "E41" "  +1.79%  "
Set-TextValue "D42" "8.971"
Set-TextValue "E42" "  +3.16%  "
Set-TextValue "E43" "  +4.24%  "
Set-TextValue "D44" "14.86"
Set-TextValue "E44" "  +2.14%  "
Set-TextValue "D45" "0.6611"
Set-TextValue "E45" "  +8.74%  "
Set-TextValue "D46" "0.9999"
Set-TextValue "E46" "  +0.59%  "
Set-TextValue "D47" "3.956"
Set-TextValue "E47" "  +1.28%  "
Set-TextValue "D48" "2.191"
Set-TextValue "E48" "  +5.81%  "
Set-TextValue "D49" "132.72"
Set-TextValue "E49" "  +2.38%  "
Set-TextValue "D50" "0.07340"
Set-TextValue "E50" "  -0.09%  "
Set-TextValue "D51" "80.70"
Set-TextValue "E51" "  +3.28%  "
